$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 24
$ws.Range("D4").Value = "1. 从solder=cur节点开始，查找需要交换的两个节点`n2. 两个节点有一个为空，就退出程序，返回solder.next`n3. 如果这两个节点全部不为空，交换两个节点，cur指针前进两步`n4. 进入步骤2"
$ws.Range("C4").Value = "给定一个链表，两两交换其中相邻的节点，并返回交换后的链表"
$ws.Range("E4").Value = "solder`n链表插入`n节点交换"
$ws.Range("F4").Value = "O(N), N是元素个数"
$ws.Range("G4").Value = "O(1)"

$ws.Rows.Item(4).RowHeight = 120

$ws.Range("D3").Select()
